{"js": "// Remove the trailing \"footer\" paragraphs that used to follow the last\n// requirement line (\"LOQ4083: Fen\u00f4menos de Transporte I (Requisito fraco)\"):\n//   - an empty paragraph\n//   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//      pages. Original theme under Creative Commons Attribution\"\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\n// Find the paragraph that holds the last requirement line; the footer block\n// we need to remove starts right after it (an empty paragraph followed by\n// the two text paragraphs above).\nlet anchorIndex = -1;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === \"LOQ4083: Fen\u00f4menos de Transporte I (Requisito fraco)\") {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex !== -1) {\n  // Delete the 3 paragraphs that immediately follow the anchor: the blank\n  // paragraph, then the two footer-text paragraphs.\n  for (let i = 0; i < 3; i++) {\n    const p = paragraphs.items[anchorIndex + 1];\n    if (!p) break;\n    p.delete();\n  }\n} else {\n  // Fallback: delete by matching text directly, plus the empty paragraph\n  // immediately preceding the \"Ver no Jupiter...\" paragraph.\n  for (let i = paragraphs.items.length - 1; i >= 0; i--) {\n    const text = paragraphs.items[i].text;\n    if (targets.indexOf(text) !== -1) {\n      paragraphs.items[i].delete();\n      if (i > 0 && paragraphs.items[i - 1].text === \"\") {\n        paragraphs.items[i - 1].delete();\n      }\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"footer\" paragraphs that used to follow the last\n# requirement line (\"LOQ4083: Fen\u00f4menos de Transporte I (Requisito fraco)\"):\n#   - an empty paragraph\n#   - \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   - \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#      pages. Original theme under Creative Commons Attribution\"\n$d = $word.ActiveDocument\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($t -eq \"LOQ4083: Fen\u00f4menos de Transporte I (Requisito fraco)\") {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -ge 1) {\n    # Delete the 3 paragraphs immediately following the anchor paragraph:\n    # the blank paragraph, then the two footer-text paragraphs. Deleting\n    # repeatedly at (anchorIndex + 1) works because each delete shifts the\n    # following paragraphs up by one.\n    for ($n = 0; $n -lt 3; $n++) {\n        if ($anchorIndex + 1 -le $d.Paragraphs.Count) {\n            $d.Paragraphs.Item($anchorIndex + 1).Range.Delete()\n        }\n    }\n}\n\n$d.Save()\n"}
